# Refresh the cryptocurrency price/volume snapshot (price + 1h change) pulled
# from coinranking.com, matching the GitHub Actions scheduled data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.568.52"
$ws.Range("E2").Value = "  +1.57%  "
$ws.Range("D3").Value = "1.828.41"
$ws.Range("E3").Value = "  +1.86%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'317.81"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D7").Value = "'0.5420"
$ws.Range("E7").Value = "  +0.49%  "
$ws.Range("D8").Value = "'0.4020"
$ws.Range("E8").Value = "  +6.29%  "
$ws.Range("D9").Value = "'0.07680"
$ws.Range("E9").Value = "  +3.09%  "
$ws.Range("D10").Value = "'1.122"
$ws.Range("E10").Value = "  +2.53%  "
$ws.Range("E11").Value = "  +0.40%  "
$ws.Range("D12").Value = "'21.29"
$ws.Range("E12").Value = "  +3.60%  "
$ws.Range("D13").Value = "'6.342"
$ws.Range("E13").Value = "  +3.71%  "
$ws.Range("D14").Value = "'7.650"
$ws.Range("E14").Value = "  +5.68%  "
$ws.Range("D15").Value = "'1.001"
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("D16").Value = "1.826.46"
$ws.Range("E16").Value = "  +2.10%  "
$ws.Range("D17").Value = "'0.00001093"
$ws.Range("E17").Value = "  +3.09%  "
$ws.Range("D18").Value = "'90.10"
$ws.Range("E18").Value = "  +1.15%  "
$ws.Range("D19").Value = "'0.06606"
$ws.Range("E19").Value = "  +1.73%  "
$ws.Range("D20").Value = "'17.85"
$ws.Range("E20").Value = "  +3.43%  "
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").Value = "'6.081"
$ws.Range("E22").Value = "  +3.07%  "
$ws.Range("D23").Value = "28.583.46"
$ws.Range("E23").Value = "  +1.56%  "
$ws.Range("D24").Value = "'11.20"
$ws.Range("E24").Value = "  +0.26%  "
$ws.Range("D25").Value = "'2.275"
$ws.Range("E25").Value = "  +8.96%  "
$ws.Range("D26").Value = "'158.44"
$ws.Range("E26").Value = "  +2.12%  "
$ws.Range("D27").Value = "'2.457"
$ws.Range("E27").Value = "  +7.68%  "
$ws.Range("D28").Value = "'20.79"
$ws.Range("E28").Value = "  +2.45%  "
$ws.Range("D29").Value = "2.037.09"
$ws.Range("E29").Value = "  +2.14%  "
$ws.Range("D30").Value = "'124.29"
$ws.Range("E30").Value = "  +2.60%  "
$ws.Range("D31").Value = "'1.133"
$ws.Range("E31").Value = "  +1.17%  "
$ws.Range("D32").Value = "'0.1114"
$ws.Range("E32").Value = "  +4.88%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.07566"
$ws.Range("E33").Value = "  +16.61%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'5.690"
$ws.Range("E34").Value = "  +2.49%  "
$ws.Range("D35").Value = "'3.647"
$ws.Range("E35").Value = "  -0.22%  "
$ws.Range("D36").Value = "'0.2251"
$ws.Range("E36").Value = "  -0.21%  "
$ws.Range("D37").Value = "'0.02362"
$ws.Range("E37").Value = "  +3.01%  "
$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").Value = "'8.912"
$ws.Range("E38").Value = "  +5.49%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "'5.216"
$ws.Range("E39").Value = "  +4.06%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "'0.6313"
$ws.Range("E40").Value = "  +2.00%  "
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").Value = "'11.38"
$ws.Range("E41").Value = "  +2.49%  "
$ws.Range("E42").Value = "  +1.03%  "
$ws.Range("D43").Value = "'0.9999"
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("E44").Value = "  -3.28%  "
$ws.Range("D45").Value = "'13.50"
$ws.Range("E45").Value = "  +1.47%  "
$ws.Range("D46").Value = "'0.5897"
$ws.Range("E46").Value = "  +1.99%  "
$ws.Range("D47").Value = "'3.712"
$ws.Range("E47").Value = "  +1.07%  "
$ws.Range("D48").Value = "'125.50"
$ws.Range("E48").Value = "  +1.01%  "
$ws.Range("D49").Value = "'2.007"
$ws.Range("E49").Value = "  +4.16%  "
$ws.Range("D50").Value = "'1.201"
$ws.Range("E50").Value = "  +0.92%  "
$ws.Range("D51").Value = "'0.06922"
$ws.Range("E51").Value = "  +1.65%  "
